$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new rows (170-171) get the same date number format as existing date column D cells
$dateFormat = $ws.Range("D156").NumberFormat

# Row 156
$ws.Range("A156").Value = 5
$ws.Range("B156").Value = "Macroferia Regional de Talca"
$ws.Range("C156").Value = "Maule"
$ws.Range("D156").Value = 44585
$ws.Range("E156").Value = 7
$ws.Range("F156").Value = 100112021
$ws.Range("G156").Value = "Ají"
$ws.Range("H156").Value = "Americana (o)"
$ws.Range("I156").Value = "Primera"
$ws.Range("J156").Value = 100
$ws.Range("K156").Value = 12000
$ws.Range("L156").Value = 12000
$ws.Range("M156").Value = 12000
$ws.Range("N156").Value = "`$/caja 15 kilos"
$ws.Range("O156").Value = "Región del Maule"
$ws.Range("P156").Value = 800
$ws.Range("Q156").Value = 15
$ws.Range("R156").Value = "Hortaliza"

# Row 157
$ws.Range("A157").Value = 5
$ws.Range("B157").Value = "Macroferia Regional de Talca"
$ws.Range("C157").Value = "Maule"
$ws.Range("D157").Value = 44585
$ws.Range("E157").Value = 7
$ws.Range("F157").Value = 100112021
$ws.Range("G157").Value = "Ají"
$ws.Range("H157").Value = "Americana (o)"
$ws.Range("I157").Value = "Primera"
$ws.Range("J157").Value = 150
$ws.Range("K157").Value = 20000
$ws.Range("L157").Value = 20000
$ws.Range("M157").Value = 20000
$ws.Range("N157").Value = "`$/saco 25 kilos"
$ws.Range("O157").Value = "Limache"
$ws.Range("P157").Value = 800
$ws.Range("Q157").Value = 25
$ws.Range("R157").Value = "Hortaliza"

# Row 158
$ws.Range("A158").Value = 5
$ws.Range("B158").Value = "Macroferia Regional de Talca"
$ws.Range("C158").Value = "Maule"
$ws.Range("D158").Value = 44560
$ws.Range("E158").Value = 7
$ws.Range("F158").Value = 100112021
$ws.Range("G158").Value = "Ají"
$ws.Range("H158").Value = "Americana (o)"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 100
$ws.Range("K158").Value = 20000
$ws.Range("L158").Value = 20000
$ws.Range("M158").Value = 20000
$ws.Range("N158").Value = "`$/caja 14 kilos"
$ws.Range("O158").Value = "Región del Maule"
$ws.Range("P158").Value = 1429
$ws.Range("Q158").Value = 14
$ws.Range("R158").Value = "Hortaliza"

# Row 159
$ws.Range("A159").Value = 5
$ws.Range("B159").Value = "Macroferia Regional de Talca"
$ws.Range("C159").Value = "Maule"
$ws.Range("D159").Value = 44272
$ws.Range("E159").Value = 7
$ws.Range("F159").Value = 100112021
$ws.Range("G159").Value = "Ají"
$ws.Range("H159").Value = "Cacho cabra verde"
$ws.Range("I159").Value = "Primera"
$ws.Range("J159").Value = 100
$ws.Range("K159").Value = 12000
$ws.Range("L159").Value = 12000
$ws.Range("M159").Value = 12000
$ws.Range("N159").Value = "`$/saco 25 kilos"
$ws.Range("O159").Value = "Región del Maule"
$ws.Range("P159").Value = 480
$ws.Range("Q159").Value = 25
$ws.Range("R159").Value = "Hortaliza"

# Row 160
$ws.Range("A160").Value = 5
$ws.Range("B160").Value = "Macroferia Regional de Talca"
$ws.Range("C160").Value = "Maule"
$ws.Range("D160").Value = 44272
$ws.Range("E160").Value = 7
$ws.Range("F160").Value = 100112021
$ws.Range("G160").Value = "Ají"
$ws.Range("H160").Value = "Cristal"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 100
$ws.Range("K160").Value = 12000
$ws.Range("L160").Value = 12000
$ws.Range("M160").Value = 12000
$ws.Range("N160").Value = "`$/saco 25 kilos"
$ws.Range("O160").Value = "Región del Maule"
$ws.Range("P160").Value = 480
$ws.Range("Q160").Value = 25
$ws.Range("R160").Value = "Hortaliza"

# Row 161
$ws.Range("A161").Value = 5
$ws.Range("B161").Value = "Macroferia Regional de Talca"
$ws.Range("C161").Value = "Maule"
$ws.Range("D161").Value = 44529
$ws.Range("E161").Value = 7
$ws.Range("F161").Value = 100112021
$ws.Range("G161").Value = "Ají"
$ws.Range("H161").Value = "Americana (o)"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 150
$ws.Range("K161").Value = 18000
$ws.Range("L161").Value = 18000
$ws.Range("M161").Value = 18000
$ws.Range("N161").Value = "`$/caja 15 kilos"
$ws.Range("O161").Value = "Región del Maule"
$ws.Range("P161").Value = 1200
$ws.Range("Q161").Value = 15
$ws.Range("R161").Value = "Hortaliza"

# Row 162
$ws.Range("A162").Value = 5
$ws.Range("B162").Value = "Macroferia Regional de Talca"
$ws.Range("C162").Value = "Maule"
$ws.Range("D162").Value = 44414
$ws.Range("E162").Value = 7
$ws.Range("F162").Value = 100112021
$ws.Range("G162").Value = "Ají"
$ws.Range("H162").Value = "Americana (o)"
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 100
$ws.Range("K162").Value = 80000
$ws.Range("L162").Value = 80000
$ws.Range("M162").Value = 80000
$ws.Range("N162").Value = "`$/caja 25 kilos"
$ws.Range("O162").Value = "Provincia del Elquí"
$ws.Range("P162").Value = 3200
$ws.Range("Q162").Value = 25
$ws.Range("R162").Value = "Hortaliza"

# Row 163
$ws.Range("A163").Value = 5
$ws.Range("B163").Value = "Macroferia Regional de Talca"
$ws.Range("C163").Value = "Maule"
$ws.Range("D163").Value = 44306
$ws.Range("E163").Value = 7
$ws.Range("F163").Value = 100112021
$ws.Range("G163").Value = "Ají"
$ws.Range("H163").Value = "Cacho cabra verde"
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 200
$ws.Range("K163").Value = 12000
$ws.Range("L163").Value = 12000
$ws.Range("M163").Value = 12000
$ws.Range("N163").Value = "`$/saco 25 kilos"
$ws.Range("O163").Value = "Región del Maule"
$ws.Range("P163").Value = 480
$ws.Range("Q163").Value = 25
$ws.Range("R163").Value = "Hortaliza"

# Row 164
$ws.Range("A164").Value = 5
$ws.Range("B164").Value = "Macroferia Regional de Talca"
$ws.Range("C164").Value = "Maule"
$ws.Range("D164").Value = 44189
$ws.Range("E164").Value = 7
$ws.Range("F164").Value = 100112021
$ws.Range("G164").Value = "Ají"
$ws.Range("H164").Value = "Americana (o)"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 100
$ws.Range("K164").Value = 15000
$ws.Range("L164").Value = 15000
$ws.Range("M164").Value = 15000
$ws.Range("N164").Value = "`$/caja 14 kilos"
$ws.Range("O164").Value = "Región del Maule"
$ws.Range("P164").Value = 1071
$ws.Range("Q164").Value = 14
$ws.Range("R164").Value = "Hortaliza"

# Row 165
$ws.Range("A165").Value = 5
$ws.Range("B165").Value = "Macroferia Regional de Talca"
$ws.Range("C165").Value = "Maule"
$ws.Range("D165").Value = 44299
$ws.Range("E165").Value = 7
$ws.Range("F165").Value = 100112021
$ws.Range("G165").Value = "Ají"
$ws.Range("H165").Value = "Cristal"
$ws.Range("I165").Value = "Primera"
$ws.Range("J165").Value = 100
$ws.Range("K165").Value = 14000
$ws.Range("L165").Value = 14000
$ws.Range("M165").Value = 14000
$ws.Range("N165").Value = "`$/saco 25 kilos"
$ws.Range("O165").Value = "Región del Maule"
$ws.Range("P165").Value = 560
$ws.Range("Q165").Value = 25
$ws.Range("R165").Value = "Hortaliza"

# Row 166
$ws.Range("A166").Value = 5
$ws.Range("B166").Value = "Macroferia Regional de Talca"
$ws.Range("C166").Value = "Maule"
$ws.Range("D166").Value = 44321
$ws.Range("E166").Value = 7
$ws.Range("F166").Value = 100112021
$ws.Range("G166").Value = "Ají"
$ws.Range("H166").Value = "Cacho cabra rojo"
$ws.Range("I166").Value = "Primera"
$ws.Range("J166").Value = 100
$ws.Range("K166").Value = 15000
$ws.Range("L166").Value = 15000
$ws.Range("M166").Value = 15000
$ws.Range("N166").Value = "`$/saco 25 kilos"
$ws.Range("O166").Value = "Región del Maule"
$ws.Range("P166").Value = 600
$ws.Range("Q166").Value = 25
$ws.Range("R166").Value = "Hortaliza"

# Row 167
$ws.Range("A167").Value = 5
$ws.Range("B167").Value = "Macroferia Regional de Talca"
$ws.Range("C167").Value = "Maule"
$ws.Range("D167").Value = 44302
$ws.Range("E167").Value = 7
$ws.Range("F167").Value = 100112021
$ws.Range("G167").Value = "Ají"
$ws.Range("H167").Value = "Cacho cabra verde"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 200
$ws.Range("K167").Value = 11000
$ws.Range("L167").Value = 11000
$ws.Range("M167").Value = 11000
$ws.Range("N167").Value = "`$/saco 25 kilos"
$ws.Range("O167").Value = "Región del Maule"
$ws.Range("P167").Value = 440
$ws.Range("Q167").Value = 25
$ws.Range("R167").Value = "Hortaliza"

# Row 168
$ws.Range("A168").Value = 5
$ws.Range("B168").Value = "Macroferia Regional de Talca"
$ws.Range("C168").Value = "Maule"
$ws.Range("D168").Value = 44209
$ws.Range("E168").Value = 7
$ws.Range("F168").Value = 100112021
$ws.Range("G168").Value = "Ají"
$ws.Range("H168").Value = "Americana (o)"
$ws.Range("I168").Value = "Primera"
$ws.Range("J168").Value = 300
$ws.Range("K168").Value = 14000
$ws.Range("L168").Value = 14000
$ws.Range("M168").Value = 14000
$ws.Range("N168").Value = "`$/caja 14 kilos"
$ws.Range("O168").Value = "Región del Maule"
$ws.Range("P168").Value = 1000
$ws.Range("Q168").Value = 14
$ws.Range("R168").Value = "Hortaliza"

# Row 169
$ws.Range("A169").Value = 5
$ws.Range("B169").Value = "Macroferia Regional de Talca"
$ws.Range("C169").Value = "Maule"
$ws.Range("D169").Value = 44274
$ws.Range("E169").Value = 7
$ws.Range("F169").Value = 100112021
$ws.Range("G169").Value = "Ají"
$ws.Range("H169").Value = "Cacho cabra verde"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 100
$ws.Range("K169").Value = 12000
$ws.Range("L169").Value = 12000
$ws.Range("M169").Value = 12000
$ws.Range("N169").Value = "`$/saco 25 kilos"
$ws.Range("O169").Value = "Región del Maule"
$ws.Range("P169").Value = 480
$ws.Range("Q169").Value = 25
$ws.Range("R169").Value = "Hortaliza"

# Row 170
$ws.Range("A170").Value = 5
$ws.Range("B170").Value = "Macroferia Regional de Talca"
$ws.Range("C170").Value = "Maule"
$ws.Range("D170").Value = 44274
$ws.Range("E170").Value = 7
$ws.Range("F170").Value = 100112021
$ws.Range("G170").Value = "Ají"
$ws.Range("H170").Value = "Cristal"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 100
$ws.Range("K170").Value = 12000
$ws.Range("L170").Value = 12000
$ws.Range("M170").Value = 12000
$ws.Range("N170").Value = "`$/saco 25 kilos"
$ws.Range("O170").Value = "Región del Maule"
$ws.Range("P170").Value = 480
$ws.Range("Q170").Value = 25
$ws.Range("R170").Value = "Hortaliza"
$ws.Range("D170").NumberFormat = $dateFormat

# Row 171
$ws.Range("A171").Value = 5
$ws.Range("B171").Value = "Macroferia Regional de Talca"
$ws.Range("C171").Value = "Maule"
$ws.Range("D171").Value = 44554
$ws.Range("E171").Value = 7
$ws.Range("F171").Value = 100112021
$ws.Range("G171").Value = "Ají"
$ws.Range("H171").Value = "Americana (o)"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 150
$ws.Range("K171").Value = 14000
$ws.Range("L171").Value = 14000
$ws.Range("M171").Value = 14000
$ws.Range("N171").Value = "`$/caja 14 kilos"
$ws.Range("O171").Value = "Región del Maule"
$ws.Range("P171").Value = 1000
$ws.Range("Q171").Value = 14
$ws.Range("R171").Value = "Hortaliza"
$ws.Range("D171").NumberFormat = $dateFormat
